$wb = $excel.ActiveWorkbook

$wsLayout = $wb.Worksheets.Item("Layout")
$wsCalcs = $wb.Worksheets.Item("Calcs")
$wsNumbers = $wb.Worksheets.Item("Numbers")

# --- Update Layout sheet grid values (color indices) ---
$wsLayout.Range("B2").Value = 1
$wsLayout.Range("C2").Value = 1
$wsLayout.Range("D2").Value = 1
$wsLayout.Range("E2").Value = 1
$wsLayout.Range("F2").Value = 1
$wsLayout.Range("G2").Value = 1
$wsLayout.Range("H2").Value = 1
$wsLayout.Range("I2").Value = 1
$wsLayout.Range("J2").Value = 1
$wsLayout.Range("K2").Value = 1
$wsLayout.Range("L2").Value = 1
$wsLayout.Range("M2").Value = 1
$wsLayout.Range("N2").Value = 1
$wsLayout.Range("A4").Value = 1
$wsLayout.Range("B4").Value = 1
$wsLayout.Range("O4").Value = 2
$wsLayout.Range("P4").Value = 2
$wsLayout.Range("Q4").Value = 1
$wsLayout.Range("A6").Value = 1
$wsLayout.Range("C6").Value = 3
$wsLayout.Range("D6").Value = 3
$wsLayout.Range("E6").Value = 3
$wsLayout.Range("F6").Value = 3
$wsLayout.Range("G6").Value = 3
$wsLayout.Range("H6").Value = 3
$wsLayout.Range("I6").Value = 3
$wsLayout.Range("J6").Value = 3
$wsLayout.Range("K6").Value = 3
$wsLayout.Range("L6").Value = 3
$wsLayout.Range("O6").Value = 2
$wsLayout.Range("Q6").Value = 1
$wsLayout.Range("A8").Value = 1
$wsLayout.Range("C8").Value = 3
$wsLayout.Range("D8").Value = 3
$wsLayout.Range("E8").Value = 3
$wsLayout.Range("F8").Value = 3
$wsLayout.Range("G8").Value = 3
$wsLayout.Range("H8").Value = 3
$wsLayout.Range("I8").Value = 3
$wsLayout.Range("J8").Value = 3
$wsLayout.Range("K8").Value = 3
$wsLayout.Range("L8").Value = 2
$wsLayout.Range("M8").Value = 2
$wsLayout.Range("N8").Value = 1
$wsLayout.Range("Q8").Value = 1
$wsLayout.Range("A10").Value = 1
$wsLayout.Range("D10").Value = 3
$wsLayout.Range("E10").Value = 3
$wsLayout.Range("F10").Value = 3
$wsLayout.Range("G10").Value = 3
$wsLayout.Range("H10").Value = 3
$wsLayout.Range("I10").Value = 3
$wsLayout.Range("J10").Value = 3
$wsLayout.Range("K10").Value = 2
$wsLayout.Range("L10").Value = 2
$wsLayout.Range("M10").Value = 2
$wsLayout.Range("N10").Value = 1
$wsLayout.Range("Q10").Value = 1
$wsLayout.Range("A12").Value = 1
$wsLayout.Range("B12").Value = 1
$wsLayout.Range("C12").Value = 1
$wsLayout.Range("D12").Value = 1
$wsLayout.Range("E12").Value = 1
$wsLayout.Range("L12").Value = 1
$wsLayout.Range("M12").Value = 1
$wsLayout.Range("Q12").Value = 1
$wsLayout.Range("B14").Value = 1
$wsLayout.Range("C14").Value = 1
$wsLayout.Range("D14").Value = 1
$wsLayout.Range("E14").Value = 1
$wsLayout.Range("F14").Value = 1
$wsLayout.Range("G14").Value = 1
$wsLayout.Range("H14").Value = 1
$wsLayout.Range("I14").Value = 1
$wsLayout.Range("J14").Value = 1
$wsLayout.Range("K14").Value = 1
$wsLayout.Range("L14").Value = 1
$wsLayout.Range("M14").Value = 1
$wsLayout.Range("N14").Value = 1
$wsLayout.Range("O14").Value = 1
$wsLayout.Range("P14").Value = 1

# --- Update sheet selections / active tab ---
$wsLayout.Select()
$wsLayout.Range("B4").Select()

$wsCalcs.Select()
$wsCalcs.Range("D6").Select()

$wsNumbers.Select()
$wsNumbers.Range("G2").Select()
